$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.438977718353271
$ws.Range("B1").Value = 1.232601284980774
$ws.Range("C1").Value = 4.612543106079102
$ws.Range("D1").Value = 2.177448272705078
$ws.Range("E1").Value = 0.7215542793273926
